$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force these D-column cells to remain Text (otherwise Excel would
# auto-convert plain numeric-looking strings like "42.19" into numbers).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply cell value updates per the diff
$ws.Range("D2").Value = '45.469.64'
$ws.Range("E2").Value = '  +6.56%  '
$ws.Range("D3").Value = '2.381.94'
$ws.Range("E3").Value = '  +4.79%  '
$ws.Range("E4").Value = '  +0.25%  '
$ws.Range("B5").Value = 'Solana'
$ws.Range("C5").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D5").Value = '111.77'
$ws.Range("E5").Value = '  +8.67%  '
$ws.Range("B6").Value = 'BNB'
$ws.Range("C6").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("D6").Value = '317.56'
$ws.Range("E6").Value = '  +2.80%  '
$ws.Range("E7").Value = '  +2.92%  '
$ws.Range("E8").Value = '  -0.38%  '
$ws.Range("D9").Value = '0.631'
$ws.Range("E9").Value = '  +5.76%  '
$ws.Range("D10").Value = '42.19'
$ws.Range("E10").Value = '  +9.97%  '
$ws.Range("D11").Value = '0.0932'
$ws.Range("E11").Value = '  +4.19%  '
$ws.Range("E12").Value = '  +6.49%  '
$ws.Range("E13").Value = '  +5.08%  '
$ws.Range("E14").Value = '  +1.37%  '
$ws.Range("D15").Value = '15.79'
$ws.Range("E15").Value = '  +5.78%  '
$ws.Range("D16").Value = '2.743.69'
$ws.Range("E16").Value = '  +4.85%  '
$ws.Range("D17").Value = '2.396.89'
$ws.Range("E17").Value = '  +5.74%  '
$ws.Range("D18").Value = '45.405.41'
$ws.Range("E18").Value = '  +6.78%  '
$ws.Range("D19").Value = '7.68'
$ws.Range("E19").Value = '  +7.10%  '
$ws.Range("E20").Value = '  +4.89%  '
$ws.Range("D21").Value = '13.13'
$ws.Range("E21").Value = '  +0.81%  '
$ws.Range("D22").Value = '75.13'
$ws.Range("E22").Value = '  +3.65%  '
$ws.Range("E23").Value = '  +5.31%  '
$ws.Range("D24").Value = '270.21'
$ws.Range("E24").Value = '  +3.60%  '
$ws.Range("D25").Value = '2.33'
$ws.Range("E25").Value = '  +8.33%  '
$ws.Range("E26").Value = '  -0.78%  '
$ws.Range("D27").Value = '11.29'
$ws.Range("E27").Value = '  +7.10%  '
$ws.Range("D28").Value = '7.53'
$ws.Range("E28").Value = '  +10.32%  '
$ws.Range("E29").Value = '  +0.44%  '
$ws.Range("D30").Value = '22.94'
$ws.Range("E30").Value = '  +4.18%  '
$ws.Range("D31").Value = '38.67'
$ws.Range("E31").Value = '  +9.46%  '
$ws.Range("E32").Value = '  +12.00%  '
$ws.Range("D33").Value = '169.93'
$ws.Range("E33").Value = '  +3.95%  '
$ws.Range("D34").Value = '3.00'
$ws.Range("E34").Value = '  +17.76%  '
$ws.Range("E35").Value = '  +3.64%  '
$ws.Range("E36").Value = '  +6.86%  '
$ws.Range("D37").Value = '4.88'
$ws.Range("E37").Value = '  +9.26%  '
$ws.Range("D38").Value = '3.07'
$ws.Range("E38").Value = '  +13.46%  '
$ws.Range("E39").Value = '  +5.96%  '
$ws.Range("E40").Value = '  +7.74%  '
$ws.Range("E41").Value = '  +13.24%  '
$ws.Range("D42").Value = '105.13'
$ws.Range("E42").Value = '  +7.26%  '
$ws.Range("B43").Value = 'Algorand'
$ws.Range("C43").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D43").Value = '0.241'
$ws.Range("E43").Value = '  +7.86%  '
$ws.Range("B44").Value = 'Celestia'
$ws.Range("C44").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D44").Value = '13.54'
$ws.Range("E44").Value = '  +15.08%  '
$ws.Range("D45").Value = '71.27'
$ws.Range("E45").Value = '  +5.12%  '
$ws.Range("E46").Value = '  -0.36%  '
$ws.Range("D47").Value = '118.63'
$ws.Range("E47").Value = '  +8.92%  '
$ws.Range("D48").Value = '5.87'
$ws.Range("E48").Value = '  +15.24%  '
$ws.Range("D49").Value = '1.66'
$ws.Range("E49").Value = '  +22.64%  '
$ws.Range("D50").Value = '9.30'
$ws.Range("E50").Value = '  +8.84%  '
$ws.Range("D51").Value = '79.02'
$ws.Range("E51").Value = '  +4.82%  '
